# Update the Metadata sheet: Experimental flag, Date, and Description value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental (row 7) was previously blank -> now the literal text "false".
# Assigning the bare text "false" directly gets auto-coerced to a Boolean by
# Excel's input parser, so instead we write it as a text-literal formula and
# then paste-special just the values back over itself to collapse the
# formula down to a plain text cell (keeps the existing cell style too).
$cellExperimental = $ws.Range("B7")
$cellExperimental.Formula = "=""false"""
$cellExperimental.Copy()
$cellExperimental.PasteSpecial(-4163)

# Date (row 8) gets refreshed to the new generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description (row 17) was previously blank -> now has the CodeSystem description.
$ws.Range("B17").Value = "Codes for comparing current values to personal baseline"
